$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column C (these_these), which becomes
# the new "these_kurzform" column. All existing columns from C onward (these_these,
# SPD, CDU / CSU, ... WerteUnion) shift one column to the right (C->D, D->E, ... AE->AF).
$ws.Columns("C:C").Insert()

# Header for the new column
$ws.Range("C1").Value = "these_kurzform"

# Short-form ("Kurzform") thesis texts, one per data row (rows 2-39)
$kurzformen = @(
    'Militärhilfe für Ukraine fortsetzen',
    'Erneuerbare Energien finanziell fördern',
    'Bürgergeld bei Ablehnung Stellenangebot streichen',
    'Tempolimit auf Autobahnen einführen',
    'aus EU-Staat eingereiste Asylsuchende abweisen',
    'Mietpreisbegrenzung beibehalten',
    'Gesichtserkennung an Bahnhöfen erlauben',
    'Stromkosten-Ausgleich für Unternehmen',
    'Rentenabschläge nach 40 Jahren abschaffen',
    '„Verantwortung vor Gott“ im GG beibehalten',
    'Fachkräfteanwerbung aus Ausland fördern',
    'Kernenergie zur Stromerzeugung nutzen',
    'Spitzensteuersatz anheben',
    'Schulpolitik-Befugnisse für Bund stärken',
    'Rüstungsexporte nach Israel erlauben',
    'Gesetzliche Krankenversicherungspflicht für alle',
    'Frauenquote in Vorständen abschaffen',
    'Ökologische Landwirtschaft stärker fördern',
    'Rechtsextremismus-Projekte verstärkt fördern',
    'Menschenrechtskontrolle durch Unternehmen',
    'BAföG einkommensabhängig belassen',
    'Schuldenbremse im Grundgesetz beibehalten',
    'Arbeitserlaubnis für Asylsuchende sofort',
    'Klimaneutralitätsziel aufgeben',
    '35-Stunden-Woche gesetzlich einführen',
    'Beratungspflicht vor Schwangerschaftsabbruch',
    'Euro durch nationale Währung ersetzen',
    'Schienenverkehr vor Straßenverkehr priorisieren',
    'Ehrenamtliche Tätigkeiten für Rente anrechnen',
    'Grundsteuer auf Mieter umlegen',
    'Streikrecht in kritischen Bereichen einschränken',
    'Volksentscheide auf Bundesebene ermöglichen',
    'Strafbarkeit ab 14 einführen',
    'EU-Zölle auf E-Autos abschaffen',
    'Doppelstaatsbürgerschaft weiterhin ermöglichen',
    'Soziales Pflichtjahr für Jugendliche einführen',
    'Fossile Heizungen weiterhin erlauben',
    'Mindestlohn auf 15 Euro erhöhen'
)

$n = $kurzformen.Length
for ($i = 0; $i -lt $n; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $kurzformen[$i]
}
